$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Samin Batra"
$ws.Range("B1").Value = "sb@snu"
$ws.Range("C1").Value = "tiger"
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 0

$ws.Range("A2").Value = "Pranjal"
$ws.Range("B2").Value = "pm@snu"
$ws.Range("C2").Value = "pranjal"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:sb@snu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:pm@snu") | Out-Null

$ws.Range("E2").Select() | Out-Null
